$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# --- Fix up E11:E13 ("bsecode") which were stored as text -> make them numeric ---
$ws.Range("E11").Value = 543237
$ws.Range("E12").Value = 526371
$ws.Range("E13").Value = 532155

# --- Append three new rows (14-16) of freshly scraped screener data ---
# Row 14: MAZDOCK
$ws.Range("A14").Value = "11/06/2024 02:39:37"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "MAZDOCK"
$ws.Range("D14").Value = "Mazagon Dock Shipbuilders Ltd"
$ws.Range("E14:E16").NumberFormat = "@"
$ws.Range("E14").Value = "543237"
$ws.Range("F14").Value = -0.92
$ws.Range("G14").Value = 3124
$ws.Range("H14").Value = 1372527

# Row 15: NMDC
$ws.Range("A15").Value = "11/06/2024 02:39:37"
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = "NMDC"
$ws.Range("D15").Value = "Nmdc Limited"
$ws.Range("E15").Value = "526371"
$ws.Range("F15").Value = -1.47
$ws.Range("G15").Value = 254.7
$ws.Range("H15").Value = 5191811

# Row 16: GAIL
$ws.Range("A16").Value = "11/06/2024 02:39:37"
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = "GAIL"
$ws.Range("D16").Value = "Gail (india) Limited"
$ws.Range("E16").Value = "532155"
$ws.Range("F16").Value = -2.13
$ws.Range("G16").Value = 208.18
$ws.Range("H16").Value = 38291847
